$d = $word.ActiveDocument

$oldText = "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$find = $d.Content.Find
$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$r = $find.Parent
$insertAt = $r.Start
$r.Delete()

$d.Range($insertAt, $insertAt).InsertAfter($newText)
